# Auto-generated edit script: updates Leve price/profit columns (H:N)
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR to match refreshed
# market data from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 371
$ws.Range("I33").Value = 342.1
$ws.Range("J33").Value = 515.5
$ws.Range("K33").Value = 342.1
$ws.Range("L33").Value = 515.5
$ws.Range("M33").Value = -113.1
$ws.Range("N33").Value = -973.5

$ws.Range("H62").Value = 27783280
$ws.Range("I62").Value = 37043704
$ws.Range("J62").Value = 2006
$ws.Range("K62").Value = 37043704
$ws.Range("L62").Value = 2006
$ws.Range("M62").Value = -37043080
$ws.Range("N62").Value = -3254

$ws.Range("H65").Value = 27783280
$ws.Range("I65").Value = 37043704
$ws.Range("J65").Value = 2006
$ws.Range("K65").Value = 185218520
$ws.Range("L65").Value = 10030
$ws.Range("M65").Value = -185215400
$ws.Range("N65").Value = -16270

$ws.Range("H88").Value = 443019.78
$ws.Range("I88").Value = 1498.25
$ws.Range("J88").Value = 619628.4
$ws.Range("K88").Value = 1498.25
$ws.Range("L88").Value = 619628.4
$ws.Range("M88").Value = -1092.25
$ws.Range("N88").Value = -620440.4

$ws.Range("H91").Value = 443019.78
$ws.Range("I91").Value = 1498.25
$ws.Range("J91").Value = 619628.4
$ws.Range("K91").Value = 1498.25
$ws.Range("L91").Value = 619628.4
$ws.Range("M91").Value = -94.25
$ws.Range("N91").Value = -622436.4

$ws.Range("H98").Value = 2641.0476
$ws.Range("I98").Value = 2694.3076
$ws.Range("J98").Value = 1948.6666
$ws.Range("K98").Value = 2694.3076
$ws.Range("L98").Value = 1948.6666
$ws.Range("M98").Value = -1196.3076
$ws.Range("N98").Value = -4944.6666

$ws.Range("H114").Value = 38000
$ws.Range("J114").Value = 38000
$ws.Range("L114").Value = 38000
$ws.Range("N114").Value = -46678

$ws.Range("H122").Value = 2641.0476
$ws.Range("I122").Value = 2694.3076
$ws.Range("J122").Value = 1948.6666
$ws.Range("K122").Value = 8082.9228
$ws.Range("L122").Value = 5845.9998
$ws.Range("M122").Value = -5632.9228
$ws.Range("N122").Value = -10745.9998

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents() | Out-Null

$ws.Range("H132").Value = 8138020.5
$ws.Range("I132").Value = 15880512
$ws.Range("J132").Value = 8404.549999999999
$ws.Range("K132").Value = 47641536
$ws.Range("L132").Value = 25213.65
$ws.Range("M132").Value = -47639006
$ws.Range("N132").Value = -30273.65

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1087.1111
$ws.Range("I2").Value = 1005.5
$ws.Range("K2").Value = 1005.5
$ws.Range("M2").Value = -892.5

$ws.Range("H45").Value = 1249.4546
$ws.Range("I45").Value = 1077.7142
$ws.Range("K45").Value = 1077.7142
$ws.Range("M45").Value = -700.7141999999999

$ws.Range("H116").Value = 1087.1111
$ws.Range("I116").Value = 1005.5
$ws.Range("K116").Value = 1005.5
$ws.Range("M116").Value = 1288.5

$ws.Range("H132").Value = 3280.8
$ws.Range("I132").Value = 2654.6365
$ws.Range("K132").Value = 7963.9095
$ws.Range("M132").Value = -5433.9095

$ws.Range("H140").Value = 31947.777
$ws.Range("J140").Value = 31947.777
$ws.Range("L140").Value = 31947.777
$ws.Range("N140").Value = -42307.777

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1087.1111
$ws.Range("I3").Value = 1005.5
$ws.Range("K3").Value = 1005.5
$ws.Range("M3").Value = -891.5

$ws.Range("H20").Value = 2227.5715
$ws.Range("I20").Value = 1779.3529
$ws.Range("J20").Value = 4132.5
$ws.Range("K20").Value = 1779.3529
$ws.Range("L20").Value = 4132.5
$ws.Range("M20").Value = -1532.3529
$ws.Range("N20").Value = -4626.5

$ws.Range("H86").Value = 2469.9583
$ws.Range("I86").Value = 2436.45
$ws.Range("J86").Value = 2637.5
$ws.Range("K86").Value = 2436.45
$ws.Range("L86").Value = 2637.5
$ws.Range("M86").Value = -1313.45
$ws.Range("N86").Value = -4883.5

$ws.Range("H89").Value = 2469.9583
$ws.Range("I89").Value = 2436.45
$ws.Range("J89").Value = 2637.5
$ws.Range("K89").Value = 12182.25
$ws.Range("L89").Value = 13187.5
$ws.Range("M89").Value = -6566.25
$ws.Range("N89").Value = -24419.5

$ws.Range("H103").Value = 10000
$ws.Range("J103").Value = 10000
$ws.Range("L103").Value = 10000
$ws.Range("N103").Value = -12344

$ws.Range("H134").Value = 5873.913
$ws.Range("I134").Value = 1057.9474
$ws.Range("J134").Value = 28749.75
$ws.Range("K134").Value = 3173.8422
$ws.Range("L134").Value = 86249.25
$ws.Range("M134").Value = -638.8422
$ws.Range("N134").Value = -91319.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 33000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents() | Out-Null

$ws.Range("H77").Value = 33000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents() | Out-Null

$ws.Range("H114").Value = 24076.77
$ws.Range("J114").Value = 24076.77
$ws.Range("L114").Value = 24076.77
$ws.Range("N114").Value = -32754.77

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 1825.75
$ws.Range("J31").Value = 3001
$ws.Range("L31").Value = 9003
$ws.Range("N31").Value = -9579

$ws.Range("H103").Value = 18920.666
$ws.Range("I103").Value = 1512.5
$ws.Range("J103").Value = 27624.75
$ws.Range("K103").Value = 4537.5
$ws.Range("L103").Value = 82874.25
$ws.Range("M103").Value = -3658.5
$ws.Range("N103").Value = -84632.25

$ws.Range("H131").Value = 14087133
$ws.Range("I131").Value = 90909400
$ws.Range("J131").Value = 3050.9
$ws.Range("K131").Value = 272728200
$ws.Range("L131").Value = 9152.700000000001
$ws.Range("M131").Value = -272723160
$ws.Range("N131").Value = -19232.7

$ws.Range("H132").Value = 982.1667
$ws.Range("I132").Value = 974.75
$ws.Range("J132").Value = 997
$ws.Range("K132").Value = 8772.75
$ws.Range("L132").Value = 8973
$ws.Range("M132").Value = -6242.75
$ws.Range("N132").Value = -14033

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 37503384
$ws.Range("I70").Value = 31253524
$ws.Range("J70").Value = 50003100
$ws.Range("K70").Value = 31253524
$ws.Range("L70").Value = 50003100
$ws.Range("M70").Value = -31253254
$ws.Range("N70").Value = -50003640

$ws.Range("H73").Value = 37503384
$ws.Range("I73").Value = 31253524
$ws.Range("J73").Value = 50003100
$ws.Range("K73").Value = 31253524
$ws.Range("L73").Value = 50003100
$ws.Range("M73").Value = -31252588
$ws.Range("N73").Value = -50004972

$ws.Range("H128").Value = 37340
$ws.Range("J128").Value = 37340
$ws.Range("L128").Value = 37340
$ws.Range("N128").Value = -47300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8137.75
$ws.Range("I7").Value = 1620.9
$ws.Range("J7").Value = 18999.166
$ws.Range("K7").Value = 1620.9
$ws.Range("L7").Value = 18999.166
$ws.Range("M7").Value = -1508.9
$ws.Range("N7").Value = -19223.166

$ws.Range("H64").Value = 10000
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10450

$ws.Range("H67").Value = 10000
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11560

$ws.Range("H122").Value = 35735044
$ws.Range("I122").Value = 50020704
$ws.Range("K122").Value = 150062112
$ws.Range("M122").Value = -150059662

$ws.Range("H123").Value = 40944.332
$ws.Range("J123").Value = 40944.332
$ws.Range("L123").Value = 40944.332
$ws.Range("N123").Value = -50744.332

$ws.Range("H126").Value = 8137.75
$ws.Range("I126").Value = 1620.9
$ws.Range("J126").Value = 18999.166
$ws.Range("K126").Value = 4862.700000000001
$ws.Range("L126").Value = 56997.49800000001
$ws.Range("M126").Value = -2392.700000000001
$ws.Range("N126").Value = -61937.49800000001

$ws.Range("H132").Value = 57715.05
$ws.Range("I132").Value = 22980.6
$ws.Range("J132").Value = 69293.2
$ws.Range("K132").Value = 68941.79999999999
$ws.Range("L132").Value = 207879.6
$ws.Range("M132").Value = -66411.79999999999
$ws.Range("N132").Value = -212939.6

$ws.Range("H136").Value = 13372.75
$ws.Range("I136").Value = 14926
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 44778
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -42228
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2000
$ws.Range("J2").Value = 2000
$ws.Range("L2").Value = 2000
$ws.Range("N2").Value = -2224

$ws.Range("H63").Value = 13102.857
$ws.Range("I63").Value = 5613
$ws.Range("K63").Value = 5613
$ws.Range("M63").Value = -4989

$ws.Range("H66").Value = 13102.857
$ws.Range("I66").Value = 5613
$ws.Range("K66").Value = 16839
$ws.Range("M66").Value = -13719

$ws.Range("H132").Value = 2674.4443
$ws.Range("I132").Value = 2509.5715
$ws.Range("K132").Value = 7528.7145
$ws.Range("M132").Value = -4998.7145

$ws.Range("H136").Value = 1333.1428
$ws.Range("I136").Value = 1391.75
$ws.Range("J136").Value = 1255
$ws.Range("K136").Value = 4175.25
$ws.Range("L136").Value = 3765
$ws.Range("M136").Value = -1625.25
$ws.Range("N136").Value = -8865
